$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 6
$ws.Range("F2").Value = -1
$ws.Range("G2").Value = -5
$ws.Range("H2").Value = 21
$ws.Range("I2").Value = 5

# Row 3
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = -2
$ws.Range("G3").Value = -4
$ws.Range("H3").Value = 32
$ws.Range("I3").Value = 5

# Row 4
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = -5
$ws.Range("G4").Value = -1
$ws.Range("H4").Value = 65
$ws.Range("I4").Value = 5

# Row 5
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = -3
$ws.Range("G5").Value = -3
$ws.Range("H5").Value = 43
$ws.Range("I5").Value = 5

# New row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = -4
$ws.Range("G6").Value = -2
$ws.Range("H6").Value = 54
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim2_1"

# Selection matches the diff's recorded cursor position
$ws.Range("I1").Select()
